# Refresh cached Universalis market-price / leve-profit figures (columns H:N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the latest scheduled pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 50: A Patch-up Place
$ws.Range("H50").Value = 75001420
$ws.Range("I50").Value = 300000000
$ws.Range("J50").Value = 1899
$ws.Range("K50").Value = 900000000
$ws.Range("L50").Value = 5697
$ws.Range("M50").Value = -899999525
$ws.Range("N50").Value = -6647

# Row 51: A Bile Business
$ws.Range("H51").Value = 6391.0835
$ws.Range("I51").Value = 2196.2
$ws.Range("J51").Value = 7495
$ws.Range("K51").Value = 2196.2
$ws.Range("L51").Value = 7495
$ws.Range("M51").Value = -1712.2
$ws.Range("N51").Value = -8463

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 5952.6816
$ws.Range("I76").Value = 5208.0835
$ws.Range("J76").Value = 6846.2
$ws.Range("K76").Value = 5208.0835
$ws.Range("L76").Value = 6846.2
$ws.Range("M76").Value = -4893.0835
$ws.Range("N76").Value = -7476.2

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 5952.6816
$ws.Range("I79").Value = 5208.0835
$ws.Range("J79").Value = 6846.2
$ws.Range("K79").Value = 5208.0835
$ws.Range("L79").Value = 6846.2
$ws.Range("M79").Value = -4116.0835
$ws.Range("N79").Value = -9030.200000000001

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 3252.5557
$ws.Range("I100").Value = 3252.5557
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3252.5557
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2711.5557

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 1822
$ws.Range("I107").Value = 1822
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1822
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 98
$ws.Range("N107").ClearContents()

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 5262.3794
$ws.Range("I112").Value = 1149.5
$ws.Range("J112").Value = 5567.037
$ws.Range("K112").Value = 3448.5
$ws.Range("L112").Value = 16701.111
$ws.Range("M112").Value = -2340.5
$ws.Range("N112").Value = -18917.111

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2867.138
$ws.Range("I132").Value = 3236.2727
$ws.Range("J132").Value = 1707
$ws.Range("K132").Value = 9708.8181
$ws.Range("L132").Value = 5121
$ws.Range("M132").Value = -7178.8181
$ws.Range("N132").Value = -10181

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 4215.0967
$ws.Range("I137").Value = 3838
$ws.Range("J137").Value = 4394.6665
$ws.Range("K137").Value = 11514
$ws.Range("L137").Value = 13183.9995
$ws.Range("M137").Value = -8964
$ws.Range("N137").Value = -18283.9995

# Row 138: All-night Crafting
$ws.Range("H138").Value = 1967.65
$ws.Range("I138").Value = 1051.4651
$ws.Range("J138").Value = 2658.8071
$ws.Range("K138").Value = 3154.3953
$ws.Range("L138").Value = 7976.4213
$ws.Range("M138").Value = 1985.6047
$ws.Range("N138").Value = -18256.4213

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 212245.78
$ws.Range("I74").Value = 201666.2
$ws.Range("J74").Value = 225470.25
$ws.Range("K74").Value = 201666.2
$ws.Range("L74").Value = 225470.25
$ws.Range("M74").Value = -200792.2
$ws.Range("N74").Value = -227218.25

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 212245.78
$ws.Range("I77").Value = 201666.2
$ws.Range("J77").Value = 225470.25
$ws.Range("K77").Value = 1008331
$ws.Range("L77").Value = 1127351.25
$ws.Range("M77").Value = -1003963
$ws.Range("N77").Value = -1136087.25

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1686.2195
$ws.Range("I132").Value = 1274.081
$ws.Range("J132").Value = 5498.5
$ws.Range("K132").Value = 3822.242999999999
$ws.Range("L132").Value = 16495.5
$ws.Range("M132").Value = -1292.242999999999
$ws.Range("N132").Value = -21555.5

$ws = $wb.Worksheets.Item("BSM")
# Row 54: Get Me to the War on Time
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 6043.3184
$ws.Range("I86").Value = 5597.615
$ws.Range("J86").Value = 6687.1113
$ws.Range("K86").Value = 5597.615
$ws.Range("L86").Value = 6687.1113
$ws.Range("M86").Value = -4474.615
$ws.Range("N86").Value = -8933.1113

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 6043.3184
$ws.Range("I89").Value = 5597.615
$ws.Range("J89").Value = 6687.1113
$ws.Range("K89").Value = 27988.075
$ws.Range("L89").Value = 33435.5565
$ws.Range("M89").Value = -22372.075
$ws.Range("N89").Value = -44667.5565

# Row 132: Always Be Prepaired
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 5356.3335
$ws.Range("I134").Value = 1334.5
$ws.Range("J134").Value = 13400
$ws.Range("K134").Value = 4003.5
$ws.Range("L134").Value = 40200
$ws.Range("M134").Value = -1468.5
$ws.Range("N134").Value = -45270

# Row 135: Axes to the Maxes
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 137: Dagger Swagger
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2179.3333
$ws.Range("I31").Value = 1585.814
$ws.Range("J31").Value = 3160.923
$ws.Range("K31").Value = 1585.814
$ws.Range("L31").Value = 3160.923
$ws.Range("M31").Value = -1290.814
$ws.Range("N31").Value = -3750.923

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2179.3333
$ws.Range("I34").Value = 1585.814
$ws.Range("J34").Value = 3160.923
$ws.Range("K34").Value = 1585.814
$ws.Range("L34").Value = 3160.923
$ws.Range("M34").Value = -1383.814
$ws.Range("N34").Value = -3564.923

# Row 86: Birch, Please
$ws.Range("H86").Value = 8731.789000000001
$ws.Range("I86").Value = 6839.231
$ws.Range("J86").Value = 12832.333
$ws.Range("K86").Value = 6839.231
$ws.Range("L86").Value = 12832.333
$ws.Range("M86").Value = -5716.231
$ws.Range("N86").Value = -15078.333

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 8731.789000000001
$ws.Range("I89").Value = 6839.231
$ws.Range("J89").Value = 12832.333
$ws.Range("K89").Value = 34196.155
$ws.Range("L89").Value = 64161.665
$ws.Range("M89").Value = -28580.155
$ws.Range("N89").Value = -75393.66500000001

# Row 108: Just Starting Out
$ws.Range("H108").Value = 48206.25
$ws.Range("I108").Value = 42499.5
$ws.Range("J108").Value = 50108.5
$ws.Range("K108").Value = 42499.5
$ws.Range("L108").Value = 50108.5
$ws.Range("M108").Value = -38659.5
$ws.Range("N108").Value = -57788.5

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 24854.049
$ws.Range("I134").Value = 30394.531
$ws.Range("J134").Value = 5154.5557
$ws.Range("K134").Value = 91183.59299999999
$ws.Range("L134").Value = 15463.6671
$ws.Range("M134").Value = -88648.59299999999
$ws.Range("N134").Value = -20533.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 1431.0667
$ws.Range("I14").Value = 1431.0667
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 4293.2001
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -4120.2001

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1579.8
$ws.Range("I68").Value = 633
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1899
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = -1088
$ws.Range("N68").Value = -10622

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1579.8
$ws.Range("I71").Value = 633
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 5697
$ws.Range("L71").Value = 27000
$ws.Range("M71").Value = -1641
$ws.Range("N71").Value = -35112

# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 6833.3335
$ws.Range("I80").Value = 6500
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 19500
$ws.Range("L80").Value = 21000
$ws.Range("M80").Value = -18564
$ws.Range("N80").Value = -22872

# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 6833.3335
$ws.Range("I83").Value = 6500
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 58500
$ws.Range("L83").Value = 63000
$ws.Range("M83").Value = -53820
$ws.Range("N83").Value = -72360

# Row 132: More Mezcal
$ws.Range("H132").Value = 2922.4
$ws.Range("I132").Value = 2683.625
$ws.Range("J132").Value = 3081.5833
$ws.Range("K132").Value = 24152.625
$ws.Range("L132").Value = 27734.2497
$ws.Range("M132").Value = -21622.625
$ws.Range("N132").Value = -32794.2497

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 6732.3125
$ws.Range("I70").Value = 6388.6665
$ws.Range("J70").Value = 7174.143
$ws.Range("K70").Value = 6388.6665
$ws.Range("L70").Value = 7174.143
$ws.Range("M70").Value = -6118.6665
$ws.Range("N70").Value = -7714.143

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 6732.3125
$ws.Range("I73").Value = 6388.6665
$ws.Range("J73").Value = 7174.143
$ws.Range("K73").Value = 6388.6665
$ws.Range("L73").Value = 7174.143
$ws.Range("M73").Value = -5452.6665
$ws.Range("N73").Value = -9046.143

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 146378.23
$ws.Range("I122").Value = 233739.62
$ws.Range("J122").Value = 6600
$ws.Range("K122").Value = 701218.86
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -698768.86
$ws.Range("N122").Value = -24700

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 3153.4783
$ws.Range("I126").Value = 2926.3157
$ws.Range("J126").Value = 4232.5
$ws.Range("K126").Value = 8778.947100000001
$ws.Range("L126").Value = 12697.5
$ws.Range("M126").Value = -6308.947100000001
$ws.Range("N126").Value = -17637.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1041.7222
$ws.Range("I22").Value = 652.5454999999999
$ws.Range("J22").Value = 1653.2858
$ws.Range("K22").Value = 652.5454999999999
$ws.Range("L22").Value = 1653.2858
$ws.Range("M22").Value = -357.5454999999999
$ws.Range("N22").Value = -2243.2858

# Row 25: A Rush on Ringbands
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10460

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1041.7222
$ws.Range("I27").Value = 652.5454999999999
$ws.Range("J27").Value = 1653.2858
$ws.Range("K27").Value = 652.5454999999999
$ws.Range("L27").Value = 1653.2858
$ws.Range("M27").Value = -545.5454999999999
$ws.Range("N27").Value = -1867.2858

# Row 50: The Birdmen of Ishgard
$ws.Range("H50").Value = 20000
$ws.Range("I50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 20000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -19363
$ws.Range("N50").ClearContents()

# Row 131: For What Was Gleaned
$ws.Range("H131").Value = 68069.7
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 68069.7
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 68069.7
$ws.Range("N131").Value = -78149.7

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 1376.125
$ws.Range("I96").Value = 999.6667
$ws.Range("J96").Value = 1602
$ws.Range("K96").Value = 999.6667
$ws.Range("L96").Value = 1602
$ws.Range("M96").Value = 373.3333
$ws.Range("N96").Value = -4348

# Row 113: A Tender Table
$ws.Range("H113").Value = 805.1951
$ws.Range("I113").Value = 657.5599999999999
$ws.Range("J113").Value = 1035.875
$ws.Range("K113").Value = 1972.68
$ws.Range("L113").Value = 3107.625
$ws.Range("M113").Value = 197.3200000000002
$ws.Range("N113").Value = -7447.625

# Row 120: A Turban for the Ages
$ws.Range("H120").Value = 47333.668
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 47333.668
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 47333.668
$ws.Range("N120").Value = -57009.668

# Row 128: Lightening Up
$ws.Range("H128").Value = 55429
$ws.Range("I128").Value = 50000
$ws.Range("J128").Value = 56786.25
$ws.Range("K128").Value = 50000
$ws.Range("L128").Value = 56786.25
$ws.Range("M128").Value = -45020
$ws.Range("N128").Value = -66746.25

# Row 129: Lifetime of Gleaning
$ws.Range("H129").Value = 59263
$ws.Range("I129").Value = 59000
$ws.Range("J129").Value = 59394.5
$ws.Range("K129").Value = 59000
$ws.Range("L129").Value = 59394.5
$ws.Range("M129").Value = -54000
$ws.Range("N129").Value = -69394.5

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 10009.333
$ws.Range("I132").Value = 10448.0625
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 31344.1875
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -28814.1875
$ws.Range("N132").Value = -24558.5

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2216
$ws.Range("I136").Value = 1441.3334
$ws.Range("J136").Value = 4346.3335
$ws.Range("K136").Value = 4324.0002
$ws.Range("L136").Value = 13039.0005
$ws.Range("M136").Value = -1774.0002
$ws.Range("N136").Value = -18139.0005
